$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 911
    $ws.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202403/JJwRjJtf1710726239594.jpeg"
}
